$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 7.569942235946655
$ws.Range("B2").Value = 0.2483266922241177
$ws.Range("C2").Value = 0.008884032567342123
$ws.Range("D2").Value = 0.0003953944015714425

$ws.Range("A3").Value = 8.050627628962198
$ws.Range("B3").Value = 0.241196467824371
$ws.Range("C3").Value = 0.008636871973673502
$ws.Range("D3").Value = 0.00006781442914425014

$ws.Range("A4").Value = 7.154925028483073
$ws.Range("B4").Value = 0.9695102346140136
$ws.Range("C4").Value = 0.005666255950927734
$ws.Range("D4").Value = 0.002421370739839173

$ws.Range("A5").Value = 5.427624861399333
$ws.Range("B5").Value = 0.2813850455954224
$ws.Range("C5").Value = 0.002858956654866537
$ws.Range("D5").Value = 0.00004942010690644381
